$d = $word.ActiveDocument

# 1. Update the date in the first line (08.08.24 -> 07.08.24)
$d.Content.Find.Execute("08.08.24", $false, $false, $false, $false, $false, $true, 1, $false, "07.08.24", 2)

# 2. Replace the paper title paragraph; this also removes the trailing <w:br/>
#    since the whole paragraph range (including the line break) is overwritten.
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Language Model Can Listen While Speaking"

# 3. Replace the first body paragraph
$p3 = $d.Paragraphs(3)
$p3.Range.Text = " המאמר שמשך את תשומת ליבי בגלל שמו הקליט. המאמר מציע ארכיטקטורה של מודל Speech Language Model או SLM שיודע להקשיב תוך כדי שהוא מדבר, כלומר מודל full duplex (מושג מתחום התקשורת). בדרך כלל ל- SLM יש שני משטר עבודה: הקשבה או דיבור, כלומר המודל או מדבר או מקשיב. המאמר מעשיר את מרחב היכולות של SLM ומצייד אותו ביכולת להקשיב תוך כדי שהוא מדבר. מעניין שהמודל גם יכול לעצור אם הוא מזהה שיש דיבור (לא רעש) ומגיב עליו (בדיבור) לאחר מכן. "

# 4. Replace the second body paragraph
$p4 = $d.Paragraphs(4)
$p4.Range.Text = "הארכיטקטורה של המודל המוצע LSLM מורכב מרכיבים סטנדרטיים. יש מודל שקולט אות דיבור, מחלק אותו לטוקנים (האות במקטעי זמן שונים) מקודד אותו לוקטור אמבדינג ומאזין אותו לדקודר. תפקיד הדקודר הוא לקחת בחשבון את ייצוג של טוקני הדיבור שנקלטו קודם וגם ייצוג טוקני הדיבור שנוצרו על ידי המודל כדי ליצור את הפלט הבא (אות הדיבור) של המודל. ֿכאמור לפעמים הדקודר מחליט שהוא צריך לעבור למצב האזנה ולפעמים הוא צריך לעבור למצב הדיבור."

# 5. Replace the third body paragraph
$p5 = $d.Paragraphs(5)
$p5.Range.Text = "כלומר הדקודר במקרה הזה הוא vocoder המקבל כקלט את אות הדיבור הנקלט בנוסף לאות הדיבור המגונרט על ה-vocoder עצמו לפני. "

# 6. Update the arxiv link
$d.Content.Find.Execute("https://arxiv.org/abs/2402.10793", $false, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/pdf/2408.02622", 2)

Write-Output "edit complete"
